$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 462, shifting the existing rows 462:475 down to 463:476
$ws.Rows(462).Insert()

# Populate the newly inserted row 462 with the new weekly price record
$ws.Cells.Item(462, 1).Value  = 5
$ws.Cells.Item(462, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(462, 3).Value  = "Maule"
$ws.Cells.Item(462, 4).Value  = 45075
$ws.Cells.Item(462, 5).Value  = 7
$ws.Cells.Item(462, 6).Value  = 100112006
$ws.Cells.Item(462, 7).Value  = "Repollo"
$ws.Cells.Item(462, 8).Value  = "Crespo record"
$ws.Cells.Item(462, 9).Value  = "Primera"
$ws.Cells.Item(462, 10).Value = 6000
$ws.Cells.Item(462, 11).Value = 700
$ws.Cells.Item(462, 12).Value = 800
$ws.Cells.Item(462, 13).Value = 750
$ws.Cells.Item(462, 14).Value = "`$/unidad"
$ws.Cells.Item(462, 15).Value = "Región del Maule"
$ws.Cells.Item(462, 16).Value = 750
$ws.Cells.Item(462, 17).Value = 1
$ws.Cells.Item(462, 18).Value = "Hortaliza"
